# Add 11 new time-log entries (rows 23-33) to Table1 on Sheet1, covering the
# work done on the iOS banking-app project from 11/9/2023 through 12/12/2023.
# Each new row is created via ListRows.Add() (so Table1's range + AutoFilter
# auto-extend), formats are copied from an existing row with a matching
# Duration style (row 9 = whole-hour "General" style, row 22 = fractional
# "h:mm" style), and then Date/Duration/Task/Description are filled in.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$rowRanges = @{}

# create row 23
$newRow = $lo.ListRows.Add()
$rr = $newRow.Range
$ws.Rows.Item(9).Copy()
$rr.PasteSpecial(-4122)
$rr.Cells.Item(1,1).Value = 45239
$rr.Cells.Item(1,2).Value = 2
$ws.Rows.Item(23).RowHeight = 45
$rowRanges[23] = $rr

# create row 24
$newRow = $lo.ListRows.Add()
$rr = $newRow.Range
$ws.Rows.Item(22).Copy()
$rr.PasteSpecial(-4122)
$rr.Cells.Item(1,1).Value = 45240
$rr.Cells.Item(1,2).Value = 0.020833333333333332
$ws.Rows.Item(24).RowHeight = 45
$rowRanges[24] = $rr

# create row 25
$newRow = $lo.ListRows.Add()
$rr = $newRow.Range
$ws.Rows.Item(9).Copy()
$rr.PasteSpecial(-4122)
$rr.Cells.Item(1,1).Value = 45246
$rr.Cells.Item(1,2).Value = 3
$ws.Rows.Item(25).RowHeight = 60
$rowRanges[25] = $rr

# create row 26
$newRow = $lo.ListRows.Add()
$rr = $newRow.Range
$ws.Rows.Item(9).Copy()
$rr.PasteSpecial(-4122)
$rr.Cells.Item(1,1).Value = 45247
$rr.Cells.Item(1,2).Value = 1
$ws.Rows.Item(26).RowHeight = 36
$rowRanges[26] = $rr

# create row 27
$newRow = $lo.ListRows.Add()
$rr = $newRow.Range
$ws.Rows.Item(9).Copy()
$rr.PasteSpecial(-4122)
$rr.Cells.Item(1,1).Value = 45256
$rr.Cells.Item(1,2).Value = 4
$ws.Rows.Item(27).RowHeight = 48
$rowRanges[27] = $rr

# create row 28
$newRow = $lo.ListRows.Add()
$rr = $newRow.Range
$ws.Rows.Item(22).Copy()
$rr.PasteSpecial(-4122)
$rr.Cells.Item(1,1).Value = 45260
$rr.Cells.Item(1,2).Value = 0.0625
$ws.Rows.Item(28).RowHeight = 36
$rowRanges[28] = $rr

# create row 29
$newRow = $lo.ListRows.Add()
$rr = $newRow.Range
$ws.Rows.Item(9).Copy()
$rr.PasteSpecial(-4122)
$rr.Cells.Item(1,1).Value = 45261
$rr.Cells.Item(1,2).Value = 3
$ws.Rows.Item(29).RowHeight = 48
$rowRanges[29] = $rr

# create row 30
$newRow = $lo.ListRows.Add()
$rr = $newRow.Range
$ws.Rows.Item(9).Copy()
$rr.PasteSpecial(-4122)
$rr.Cells.Item(1,1).Value = 45264
$rr.Cells.Item(1,2).Value = 4
$ws.Rows.Item(30).RowHeight = 36
$rowRanges[30] = $rr

# create row 31
$newRow = $lo.ListRows.Add()
$rr = $newRow.Range
$ws.Rows.Item(9).Copy()
$rr.PasteSpecial(-4122)
$rr.Cells.Item(1,1).Value = 45265
$rr.Cells.Item(1,2).Value = 1
$ws.Rows.Item(31).RowHeight = 36
$rowRanges[31] = $rr

# create row 32
$newRow = $lo.ListRows.Add()
$rr = $newRow.Range
$ws.Rows.Item(22).Copy()
$rr.PasteSpecial(-4122)
$rr.Cells.Item(1,1).Value = 45270
$rr.Cells.Item(1,2).Value = 0.013888888888888888
$ws.Rows.Item(32).RowHeight = 30
$rowRanges[32] = $rr

# create row 33
$newRow = $lo.ListRows.Add()
$rr = $newRow.Range
$ws.Rows.Item(9).Copy()
$rr.PasteSpecial(-4122)
$rr.Cells.Item(1,1).Value = 45272
$rr.Cells.Item(1,2).Value = 1
$ws.Rows.Item(33).RowHeight = 36
$rowRanges[33] = $rr

# Block 1 (rows 23-31): set all Task values first, then all Description values (matches source paste order)
$rowRanges[23].Cells.Item(1,3).Value = "Researched how to download the app"
$rowRanges[24].Cells.Item(1,3).Value = "Setting up Apple developer account"
$rowRanges[25].Cells.Item(1,3).Value = "Apple developer certificates"
$rowRanges[26].Cells.Item(1,3).Value = "Bug fixes"
$rowRanges[27].Cells.Item(1,3).Value = "user input"
$rowRanges[28].Cells.Item(1,3).Value = "TextField formatting"
$rowRanges[29].Cells.Item(1,3).Value = "Research on photo check deposit"
$rowRanges[30].Cells.Item(1,3).Value = "Research on API's"
$rowRanges[31].Cells.Item(1,3).Value = "Attempt at implementation "

$rowRanges[23].Cells.Item(1,4).Value = "watched some videos on YouTube on how to download an app onto a device from Xcode, all mentioned an apple developer account"
$rowRanges[24].Cells.Item(1,4).Value = "Set up the apple developer account and paid the yearly fee for the distribution license, this license is also needed for any development tools."
$rowRanges[25].Cells.Item(1,4).Value = "Had to make a certificate signing request from they `"Keychain`" on the mac. Then had to upload that file that was generated on apple developer website and generate a developer certificate to distribute the app"
$rowRanges[26].Cells.Item(1,4).Value = "fixed the color and other minor things after playing with the app on the phone and seeing what it really looked like"
$rowRanges[27].Cells.Item(1,4).Value = "Before even attempting coding I watched a lot of videos and did some reading on how to add user input in Xcode and what would be the best way to go about it for the type of input I needed"
$rowRanges[28].Cells.Item(1,4).Value = "After adding the TextField and had user input set up, I needed to format the the input to make it adhear to the type of number being inputted"
$rowRanges[29].Cells.Item(1,4).Value = "Adding a feature where it just opens the camera isn't hard at all, just like many things it is built into xcode. Onlt thing is I would have to get another certificate from apple to access the camera "
$rowRanges[30].Cells.Item(1,4).Value = "It looks like the easiest way to imlement a camera detecting handwritten numbers would be with some 3rd party api."
$rowRanges[31].Cells.Item(1,4).Value = "After seeing what I would need I realized not only was I not capable skill wise but also I was running out of time "

# Block 2 (rows 32-33): set Task then Description per row
$rowRanges[32].Cells.Item(1,3).Value = "Reinstalled and tested again"
$rowRanges[32].Cells.Item(1,4).Value = "Chaning the color seemed to have fixed the bug where certain elemetnts fo the app were gray"
$rowRanges[33].Cells.Item(1,3).Value = "Reviewd the code"
$rowRanges[33].Cells.Item(1,4).Value = "Went over all the pages and views and reviewd every line to make sure my code didn't work by accident and was neat and concise."

$ws.Range("A34").Select()
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1